# Update "想去人数" (number of people wanting to attend) counts on the
# "展览" and "全部类型" sheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): rows 2, 5, 6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3310
$ws1.Range("F5").Value = 1328
$ws1.Range("F6").Value = 315

# Sheet "全部类型" (4th sheet): rows 2, 5, 7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3310
$ws4.Range("F5").Value = 1328
$ws4.Range("F7").Value = 315
